$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the first three problems' "No." cell with a status fill color
# (green = solved well, red = needs review, blue = partially understood)
$ws.Range("A2").Interior.Color = 5296274   # green FF92D050
$ws.Range("A3").Interior.Color = 255       # red   FFFF0000
$ws.Range("A4").Interior.Color = 5296274   # green FF92D050
$ws.Range("A5").Interior.Color = 12611584  # blue  FF0070C0

# Add new row 5: "Merge Sorted Array" (LeetCode #88)
$ws.Range("A5").Value = 88
$ws.Range("B5").Value = "Merge Sorted Array"
$ws.Range("C5").Value = "Array"
$ws.Range("D5").Value = "Merge sort"
$ws.Range("E5").Value = "Easy"
$ws.Range("F5").Value = 1
$ws.Range("G2").Copy($ws.Range("G5"))
$ws.Range("G5").Value = "✅"
$ws.Range("H5").Value = "Given O(nlogn) sol and didn't understood the merge sort sol"

# Add understanding-percentage legend near the bottom of the sheet
$ws.Range("H18").Interior.Color = 255       # red   FFFF0000
$ws.Range("I18").Value = "0% understanding"
$ws.Range("H19").Interior.Color = 12611584  # blue  FF0070C0
$ws.Range("I19").Value = "50% understanding"
$ws.Range("H20").Interior.Color = 5296274   # green FF92D050
$ws.Range("I20").Value = "100% understanding"

$ws.Range("I21").Select()
